$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 305.55554
$ws.Range("I6").Value = 254
$ws.Range("K6").Value = 762
$ws.Range("M6").Value = -650

$ws.Range("H40").Value = 3175.6667
$ws.Range("I40").Value = 1966.4
$ws.Range("J40").Value = 5594.2
$ws.Range("K40").Value = 1966.4
$ws.Range("L40").Value = 5594.2
$ws.Range("M40").Value = -1791.4
$ws.Range("N40").Value = -5944.2

$ws.Range("H70").Value = 4130
$ws.Range("J70").Value = 4268.421
$ws.Range("L70").Value = 12805.263
$ws.Range("N70").Value = -13345.263

$ws.Range("H73").Value = 4130
$ws.Range("J73").Value = 4268.421
$ws.Range("L73").Value = 12805.263
$ws.Range("N73").Value = -14677.263

$ws.Range("H86").Value = 3171.8572
$ws.Range("I86").Value = 2533.8333
$ws.Range("J86").Value = 7000
$ws.Range("K86").Value = 2533.8333
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = -1410.8333
$ws.Range("N86").Value = -9246

$ws.Range("H89").Value = 3171.8572
$ws.Range("I89").Value = 2533.8333
$ws.Range("J89").Value = 7000
$ws.Range("K89").Value = 12669.1665
$ws.Range("L89").Value = 35000
$ws.Range("M89").Value = -7053.166499999999
$ws.Range("N89").Value = -46232

$ws.Range("H107").Value = 167.8421
$ws.Range("I107").Value = 163.5625
$ws.Range("J107").Value = 190.66667
$ws.Range("K107").Value = 163.5625
$ws.Range("L107").Value = 190.66667
$ws.Range("M107").Value = 1756.4375
$ws.Range("N107").Value = -4030.66667

$ws.Range("H111").Value = 4669.316
$ws.Range("I111").Value = 3522.0667
$ws.Range("K111").Value = 10566.2001
$ws.Range("M111").Value = -7499.2001

$ws.Range("H116").Value = 2699
$ws.Range("I116").Value = 2699
$ws.Range("K116").Value = 2699
$ws.Range("M116").Value = 743

$ws.Range("H135").Value = 5207.727
$ws.Range("I135").Value = 6128.222
$ws.Range("K135").Value = 55153.998
$ws.Range("M135").Value = -52618.998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11669.71
$ws.Range("I32").Value = 10735.704
$ws.Range("K32").Value = 10735.704
$ws.Range("M32").Value = -10448.704

$ws.Range("H110").Value = 1274.875
$ws.Range("I110").Value = 1274.875
$ws.Range("K110").Value = 1274.875
$ws.Range("M110").Value = 770.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3971.8333
$ws.Range("I86").Value = 1783.5385
$ws.Range("K86").Value = 1783.5385
$ws.Range("M86").Value = -660.5385000000001

$ws.Range("H89").Value = 3971.8333
$ws.Range("I89").Value = 1783.5385
$ws.Range("K89").Value = 8917.692500000001
$ws.Range("M89").Value = -3301.692500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 732.9167
$ws.Range("I16").Value = 545.2857
$ws.Range("J16").Value = 995.6
$ws.Range("K16").Value = 545.2857
$ws.Range("L16").Value = 995.6
$ws.Range("M16").Value = -258.2857
$ws.Range("N16").Value = -1569.6

$ws.Range("H107").Value = 4618.04
$ws.Range("I107").Value = 626.8182
$ws.Range("K107").Value = 626.8182
$ws.Range("M107").Value = 1293.1818

$ws.Range("H113").Value = 732.9167
$ws.Range("I113").Value = 545.2857
$ws.Range("J113").Value = 995.6
$ws.Range("K113").Value = 545.2857
$ws.Range("L113").Value = 995.6
$ws.Range("M113").Value = 1624.7143
$ws.Range("N113").Value = -5335.6

$ws.Range("H120").Value = 681873.75
$ws.Range("J120").Value = 681873.75
$ws.Range("L120").Value = 681873.75
$ws.Range("N120").Value = -689131.75

$ws.Range("H125").Value = 132499.5
$ws.Range("J125").Value = 132499.5
$ws.Range("L125").Value = 132499.5
$ws.Range("N125").Value = -137419.5

$ws.Range("H131").Value = 69996.5
$ws.Range("J131").Value = 69996.5
$ws.Range("L131").Value = 69996.5
$ws.Range("N131").Value = -80076.5

$ws.Range("H134").Value = 3302.392
$ws.Range("I134").Value = 1948.1666
$ws.Range("K134").Value = 5844.4998
$ws.Range("M134").Value = -3309.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 252.75
$ws.Range("I16").Value = 252.75
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 758.25
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -585.25
$ws.Range("N16").ClearContents()

$ws.Range("H37").Value = 387351.5
$ws.Range("J37").Value = 387351.5
$ws.Range("L37").Value = 1162054.5
$ws.Range("N37").Value = -1162278.5

$ws.Range("H113").Value = 8955
$ws.Range("I113").Value = 51864
$ws.Range("J113").Value = 1153.3636
$ws.Range("K113").Value = 155592
$ws.Range("L113").Value = 3460.0908
$ws.Range("M113").Value = -153422
$ws.Range("N113").Value = -7800.0908

$ws.Range("H114").Value = 1102
$ws.Range("I114").Value = 393
$ws.Range("J114").Value = 1633.75
$ws.Range("K114").Value = 1179
$ws.Range("L114").Value = 4901.25
$ws.Range("M114").Value = 2075
$ws.Range("N114").Value = -11409.25

$ws.Range("H132").Value = 1472.7273
$ws.Range("J132").Value = 1472.7273
$ws.Range("L132").Value = 13254.5457
$ws.Range("N132").Value = -18314.5457

$ws.Range("H134").Value = 8835
$ws.Range("I134").Value = 1560
$ws.Range("K134").Value = 4680
$ws.Range("M134").Value = 390

$ws.Range("H136").Value = 5309.091
$ws.Range("I136").Value = 3900
$ws.Range("K136").Value = 11700
$ws.Range("M136").Value = -6600

$ws.Range("H139").Value = 3469.1052
$ws.Range("I139").Value = 2994.2
$ws.Range("J139").Value = 5250
$ws.Range("K139").Value = 8982.599999999999
$ws.Range("L139").Value = 15750
$ws.Range("M139").Value = -3842.599999999999
$ws.Range("N139").Value = -26030

$ws.Range("H140").Value = 1616.25
$ws.Range("I140").Value = 1616.25
$ws.Range("K140").Value = 4848.75
$ws.Range("M140").Value = 331.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 149638.62
$ws.Range("I70").Value = 229422.2
$ws.Range("J70").Value = 16666
$ws.Range("K70").Value = 229422.2
$ws.Range("L70").Value = 16666
$ws.Range("M70").Value = -229152.2
$ws.Range("N70").Value = -17206

$ws.Range("H73").Value = 149638.62
$ws.Range("I73").Value = 229422.2
$ws.Range("J73").Value = 16666
$ws.Range("K73").Value = 229422.2
$ws.Range("L73").Value = 16666
$ws.Range("M73").Value = -228486.2
$ws.Range("N73").Value = -18538

$ws.Range("H107").Value = 566.4
$ws.Range("I107").Value = 341.14285
$ws.Range("K107").Value = 341.14285
$ws.Range("M107").Value = 1578.85715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 58509.5
$ws.Range("J36").Value = 58509.5
$ws.Range("L36").Value = 58509.5
$ws.Range("N36").Value = -59633.5

$ws.Range("H132").Value = 5706.154
$ws.Range("I132").Value = 2786
$ws.Range("K132").Value = 8358
$ws.Range("M132").Value = -5828

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5660.9165
$ws.Range("I81").Value = 7159.5557
$ws.Range("J81").Value = 1165
$ws.Range("K81").Value = 14319.1114
$ws.Range("L81").Value = 2330
$ws.Range("M81").Value = -13258.1114
$ws.Range("N81").Value = -4452

$ws.Range("H84").Value = 5660.9165
$ws.Range("I84").Value = 7159.5557
$ws.Range("J84").Value = 1165
$ws.Range("K84").Value = 71595.557
$ws.Range("L84").Value = 11650
$ws.Range("M84").Value = -66291.557
$ws.Range("N84").Value = -22258

$ws.Range("H107").Value = 150
$ws.Range("I107").Value = 150
$ws.Range("K107").Value = 450
$ws.Range("M107").Value = 1470

$ws.Range("H122").Value = 2801.182
$ws.Range("I122").Value = 2371.2122
$ws.Range("J122").Value = 4091.0908
$ws.Range("K122").Value = 7113.6366
$ws.Range("L122").Value = 12273.2724
$ws.Range("M122").Value = -4663.6366
$ws.Range("N122").Value = -17173.2724

$ws.Range("H132").Value = 1710.8064
$ws.Range("I132").Value = 1393.125
$ws.Range("K132").Value = 4179.375
$ws.Range("M132").Value = -1649.375

$ws.Range("H137").Value = 69272.73
$ws.Range("J137").Value = 69272.73
$ws.Range("L137").Value = 69272.73
$ws.Range("N137").Value = -79472.73

$ws.Range("H140").Value = 91578
$ws.Range("J140").Value = 91578
$ws.Range("L140").Value = 91578
$ws.Range("N140").Value = -101938

$ws.Range("H141").Value = 99686.25
$ws.Range("J141").Value = 99686.25
$ws.Range("L141").Value = 99686.25
$ws.Range("N141").Value = -110046.25
